# Tarea 06 / Vigas.xlsx — apply the commit's changes to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$CENTER = -4108   # xlCenter

# --- Row 18: add M18 ($\phi Mn [tonf \cdot m]$) with O18's style; change P18's text+style to match ---
$ws.Range("O18").Copy($ws.Range("M18")) | Out-Null
$ws.Range("M18").Value2 = "$\phi Mn [tonf \cdot m]$"

$ws.Range("O18").Copy($ws.Range("P18")) | Out-Null
$ws.Range("P18").Value2 = "$\phi \cdot Vn [tonf]$"

# --- Row 19: add M19 (number), N19 (formula, percent), P19 becomes a number, Q19 (formula, percent) ---
$ws.Range("O36").Copy($ws.Range("M19")) | Out-Null
$ws.Range("M19").Value2 = 23.67

$ws.Range("N19").Formula = "=(M19-F19)/F19"
$ws.Range("N19").NumberFormat = "0%"
$ws.Range("N19").HorizontalAlignment = $CENTER
$ws.Range("N19").VerticalAlignment = $CENTER

$ws.Range("P19").Value2 = 26.65

$ws.Range("Q19").Formula = "=(P19-F21)/F21"
$ws.Range("Q19").NumberFormat = "0%"
$ws.Range("Q19").HorizontalAlignment = $CENTER
$ws.Range("Q19").VerticalAlignment = $CENTER

# --- Row 21: add M21 ($\phi Mn [tonf \cdot m]$) with the same style as M18 ---
$ws.Range("O18").Copy($ws.Range("M21")) | Out-Null
$ws.Range("M21").Value2 = "$\phi Mn [tonf \cdot m]$"

# --- Row 22: add M22 (number), N22 (formula, percent) ---
$ws.Range("O36").Copy($ws.Range("M22")) | Out-Null
$ws.Range("M22").Value2 = 37.43

$ws.Range("N22").Formula = "=(-F20-M22)/F20"
$ws.Range("N22").NumberFormat = "0%"
$ws.Range("N22").HorizontalAlignment = $CENTER
$ws.Range("N22").VerticalAlignment = $CENTER

# --- Row 36: add Q36 (formula, percent) ---
$ws.Range("Q36").Formula = "=(O36-F36)/F36"
$ws.Range("Q36").NumberFormat = "0%"
$ws.Range("Q36").HorizontalAlignment = $CENTER
$ws.Range("Q36").VerticalAlignment = $CENTER

# --- Row 39: add Q39 (formula, percent) ---
$ws.Range("Q39").Formula = "=(O39-F37)/F37"
$ws.Range("Q39").NumberFormat = "0%"
$ws.Range("Q39").HorizontalAlignment = $CENTER
$ws.Range("Q39").VerticalAlignment = $CENTER

# --- Row 51: add Q51 (formula, percent) ---
$ws.Range("Q51").Formula = "=(O51-F54)/F54"
$ws.Range("Q51").NumberFormat = "0%"
$ws.Range("Q51").HorizontalAlignment = $CENTER
$ws.Range("Q51").VerticalAlignment = $CENTER

# --- Row 53: P53 keeps its style but loses its text ---
$ws.Range("P53").ClearContents() | Out-Null

# --- Row 54: P54 keeps its style but loses its text; add Q54 (formula, percent) ---
$ws.Range("P54").ClearContents() | Out-Null

$ws.Range("Q54").Formula = "=(O54-F53)/F53"
$ws.Range("Q54").NumberFormat = "0%"
$ws.Range("Q54").HorizontalAlignment = $CENTER
$ws.Range("Q54").VerticalAlignment = $CENTER

# --- Column Q width (new, auto-sized to content) ---
$ws.Columns.Item(17).ColumnWidth = 9.85546875

# --- Selection / active cell matches author's final view ---
$ws.Range("O19").Select() | Out-Null
